$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Primary (default) footer -> footer2.xml (Pearson logo): image2.png -> image1.png
$fPrimary = $sec.Footers.Item(1)
if ($fPrimary.Range.InlineShapes.Count -ge 1) {
    $fPrimary.Range.InlineShapes.Item(1).Name = "image1.png"
}

# First-page footer -> footer1.xml (Pearson logo): image2.png -> image1.png
$fFirst = $sec.Footers.Item(2)
if ($fFirst.Range.InlineShapes.Count -ge 1) {
    $fFirst.Range.InlineShapes.Item(1).Name = "image1.png"
}

# First-page header -> header1.xml (BTec logo): image1.jpg -> image2.jpg
$hFirst = $sec.Headers.Item(2)
if ($hFirst.Range.InlineShapes.Count -ge 1) {
    $hFirst.Range.InlineShapes.Item(1).Name = "image2.jpg"
}
